$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "threshold function" header labels.
# Written in this order so the new shared-string table comes out as
# 3/log(n), 4/log(n), 2/log(n), 1/log(n) (matching the source order the
# columns were authored in: J, M, G, D).
$ws.Range("J1").Value = "3/log(n)"
$ws.Range("M1").Value = "4/log(n)"
$ws.Range("G1").Value = "2/log(n)"
$ws.Range("D1").Value = "1/log(n)"

# Column D: 1.8/log2(n)
$ws.Range("D2").Formula = "=1.8/LOG(A2,2)"
$ws.Range("D3:D13").Formula = "=1.8/LOG(A3,2)"

# Column G: 2/log2(n)
$ws.Range("G2").Formula = "=2/LOG(A2,2)"
$ws.Range("G3:G13").Formula = "=2/LOG(A3,2)"

# Column J: 3/log2(n)
$ws.Range("J2").Formula = "=3/LOG(A2,2)"
$ws.Range("J3:J13").Formula = "=3/LOG(A3,2)"

# Column M: 4/log2(n)
$ws.Range("M2").Formula = "=4/LOG(A2,2)"
$ws.Range("M3:M13").Formula = "=4/LOG(A3,2)"

# Match the author's final cursor position.
[void]$ws.Range("J2").Select()
